# Apply the "nn res for opponent and teammate data" edit.
#
# This fills in the previously-blank Neural Network / MLP results for the
# "Opponent Data" (Optimized, columns AM:AP) and "Teammate Data" (Optimized,
# columns AU:AX) blocks on the modelIterations sheet (rows 25-27), and
# mirrors the Teammate Data numbers into the small summary table on Sheet4.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("modelIterations")
$ws4 = $wb.Worksheets.Item("Sheet4")

# --- modelIterations: Opponent Data (Optimized) + Teammate Data (Optimized) ---
# Row 25 (RMSE)
$ws1.Range("AM25").Value = 53.2
$ws1.Range("AN25").Value = 40.8
$ws1.Range("AO25").Value = 35.2
$ws1.Range("AP25").Value = 22.6
$ws1.Range("AU25").Value = 54.6
$ws1.Range("AV25").Value = 40
$ws1.Range("AW25").Value = 38.9
$ws1.Range("AX25").Value = 23.1

# Row 26 (MAE)
$ws1.Range("AM26").Value = 69.6
$ws1.Range("AN26").Value = 54.5
$ws1.Range("AO26").Value = 49.7
$ws1.Range("AP26").Value = 31.9
$ws1.Range("AU26").Value = 69.8
$ws1.Range("AV26").Value = 56.4
$ws1.Range("AW26").Value = 52.4
$ws1.Range("AX26").Value = 31.1

# Row 27 (R^2)
$ws1.Range("AM27").Value = 0.46
$ws1.Range("AN27").Value = 0.37
$ws1.Range("AO27").Value = 0.35
$ws1.Range("AP27").Value = 0.45
$ws1.Range("AU27").Value = 0.52
$ws1.Range("AV27").Value = 0.4
$ws1.Range("AW27").Value = 0.36
$ws1.Range("AX27").Value = 0.32

# --- Sheet4: summary table mirrors the Teammate Data (Optimized) MLP row ---
$ws4.Range("B1").Value = 54.6
$ws4.Range("C1").Value = 40
$ws4.Range("D1").Value = 38.9
$ws4.Range("E1").Value = 23.1

$ws4.Range("B2").Value = 69.8
$ws4.Range("C2").Value = 56.4
$ws4.Range("D2").Value = 52.4
$ws4.Range("E2").Value = 31.1

$ws4.Range("B3").Value = 0.52
$ws4.Range("C3").Value = 0.4
$ws4.Range("D3").Value = 0.36
$ws4.Range("E3").Value = 0.32

# Sheet4's banded-row fill swaps: rows 1 & 3 become white, row 2 becomes the
# light-gray shade that rows 1 & 3 used to have.
$ws4.Range("A1:E1").Interior.Color = 16777215
$ws4.Range("A3:E3").Interior.Color = 16777215
$ws4.Range("A2:E2").Interior.Color = 16119285

# --- View/selection bookkeeping ---
# Sheet4 had a lingering B1:E3 selection; clear it back to the default A1.
$ws4.Range("A1").Select()

# modelIterations stays the active sheet/tab, scrolled down a bit with a new
# active cell selection.
$ws1.Activate()
$excel.ActiveWindow.ScrollColumn = 14
$excel.ActiveWindow.ScrollRow = 10
$ws1.Range("AH28").Select()
